# --- Setup sheet references ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Warmup Plan"
$ws2 = $wb.Worksheets.Item(2)   # "Custom Domain Group"

# --- Fix font/style of Custom Domain Group!A1 (was a "code" style, now matches the row) ---
$ws2.Range("A1").Font.Name = "Calibri"
$ws2.Range("A1").Font.Color = 0

# --- Add the new "OOTB Domain Groups" sheet, after "Custom Domain Group" ---
$newws = $wb.Sheets.Add($null, $ws2)
$newws.Name = "OOTB Domain Groups"

# Phase 1: fill column A (group names) for rows 1-14; fill column B too for 2-cell rows
# (matches the order the data was actually typed in, so new shared strings come out in the same order)
$newws.Cells.Item(1, 1).Value = "Gmail"
$newws.Cells.Item(2, 1).Value = "Microsoft"
$newws.Cells.Item(3, 1).Value = "Yahoo"
$newws.Cells.Item(4, 1).Value = "Apple"
$newws.Cells.Item(5, 1).Value = "Comcast"
$newws.Cells.Item(5, 2).Value = "comcast.net"
$newws.Cells.Item(6, 1).Value = "Orange"
$newws.Cells.Item(7, 1).Value = "La Poste"
$newws.Cells.Item(7, 2).Value = "laposte.net"
$newws.Cells.Item(8, 1).Value = "Italia Online"
$newws.Cells.Item(9, 1).Value = "WP"
$newws.Cells.Item(10, 1).Value = "United Internet"
$newws.Cells.Item(11, 1).Value = "Bigpond"
$newws.Cells.Item(12, 1).Value = "Docomo"
$newws.Cells.Item(12, 2).Value = "docomo.ne.jp"
$newws.Cells.Item(13, 1).Value = "Softbank"
$newws.Cells.Item(14, 1).Value = "KDDI"

# Phase 2: fill the remaining domain lists
$row = @("gmail.com", "googlemail.com", "googlemail.co.uk")
$arr = New-Object 'object[,]' 1,3
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(1,2), $newws.Cells.Item(1,4)).Value = $arr
$row = @("live.com", "msn.com", "hotmail.ca", "hotmail.com", "hotmail.de", "hotmail.dk", "hotmail.co.jp", "hotmail.it", "hotmail.es", "hotmail.fr", "hotmail.co.uk", "hotmail.co.kr", "hotmail.com.au", "hotmail.com.ar", "hotmail.co.il", "hotmail.com.br", "hotmail.com.tr", "hotmail.co.th", "hotmail.jp", "hotmail.se", "live.at", "live.be", "live.ca", "live.cl", "live.cn", "live.co.kr", "live.co.uk", "live.com.ar", "live.com.au", "live.com.mx", "live.com.my", "live.com.sg", "live.de", "live.dk", "live.fr", "live.hk", "live.ie", "live.in", "live.it", "live.jp", "live.nl", "live.no", "live.ru", "live.se", "outlook.com", "live.co.uk", "hotmail.gr", "windowslive.com", "xbox.com", "hotmail.cl", "live.at", "live.jp", "live.ca", "hotmail.ca", "hotmail.ch", "live.fr", "live.it", "live.nl", "outlook.ie", "outlook.com.br", "live.com.pt", "live.be", "live.co.za", "mts.net")
$arr = New-Object 'object[,]' 1,64
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(2,2), $newws.Cells.Item(2,65)).Value = $arr
$row = @("yahoo.com", "rocketmail.com", "rogers.com", "sky.com", "talk21.com", "y7mail.com", "yahoo.at", "yahoo.be", "yahoo.bg", "yahoo.ca", "yahoo.cl", "yahoo.co.hu", "yahoo.co.id", "yahoo.co.il", "yahoo.co.in", "yahoo.co.jp", "yahoo.co.kr", "yahoo.com.ar", "yahoo.com.au", "yahoo.com.biz", "yahoo.com.br", "yahoo.com.cn", "yahoo.com.co", "yahoo.com.hk", "yahoo.com.hr", "yahoo.com.in", "yahoo.com.jp", "yahoo.com.kr", "yahoo.com.mx", "yahoo.com.my", "yahoo.com.net", "yahoo.com.pe", "yahoo.com.ph", "yahoo.com.sg", "yahoo.com.tr", "yahoo.com.tw", "yahoo.com.ua", "yahoo.com.ve", "yahoo.com.vn", "yahoo.co.nz", "yahoo.co.th", "yahoo.co.uk", "yahoo.co.za", "yahoo.cz", "yahoo.de", "yahoo.dk", "yahoo.ee", "yahoo.es", "yahoo.fi", "yahoo.fr", "yahoogroups.co.kr", "yahoogroups.com.cn", "yahoogroups.com.sg", "yahoogroups.com.tw", "yahoogrupper.dk", "yahoogruppi.it", "yahoo.gr", "yahoo.hr", "yahoo.hu", "yahoo.ie", "yahoo.in", "yahoo.it", "yahoo.lt", "yahoo.lv", "yahoo.nl", "yahoo.no", "yahoo.pl", "yahoo.pt", "yahoo.ro", "yahoo.rs", "yahoo.se", "yahoo.si", "yahoo.sk", "yahooxtra.co.nz", "ymail.com", "aol.com", "aim.com", "compuserve.com", "cs.com", "netscape.com", "netscape.net", "wmconnect.com", "aol.co.uk", "aol.in", "aol.de", "aol.fr", "aol.nl", "aol.pl", "aol.jp", "aol.es", "aol.it", "aol.com.ar", "aol.fi", "aol.cl", "aol.com.co", "aol.com.ve", "aol.com.au", "aol.at", "aol.be", "aol.com.br", "aol.cz", "aol.dk", "myaol.jp", "aolnorge.no", "aolpolska.pl", "aolpolcka.pl", "aolpoland.pl", "aol.ru", "aol.kr", "aol.se", "aol.ch", "aol.com.tr", "aol.co.nz", "aolchina.com", "aol.hk", "aol.tw", "luckymail.com", "verizon.net", "aol.com.mx", "bellatlantic.net", "citlink.net", "frontier.com", "frontiernet.net", "games.com", "goowy.com", "gte.net", "love.com", "verizon.net.in", "wild4music.com", "wow.com", "yahoo.cn", "yahoo.ne.jp", "yahoogroups.ca", "yahoogroups.co.in", "yahoogroups.co.uk", "yahoogroups.com", "yahoogroups.com.au", "yahoogroups.com.hk", "yahoogroups.de", "ybb.ne.jp", "ygm.com")
$arr = New-Object 'object[,]' 1,141
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(3,2), $newws.Cells.Item(3,142)).Value = $arr
$row = @("mac.com", "icloud.com", "apple.com", "me.com")
$arr = New-Object 'object[,]' 1,4
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(4,2), $newws.Cells.Item(4,5)).Value = $arr
$row = @("orange.fr", "orange.com", "wanadoo.fr", "francetelecom.com", "voila.fr", "voila.com")
$arr = New-Object 'object[,]' 1,6
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(6,2), $newws.Cells.Item(6,7)).Value = $arr
$row = @("libero.it", "inwind.it", "iol.it", "blu.it", "giallo.it", "virgilio.it")
$arr = New-Object 'object[,]' 1,6
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(8,2), $newws.Cells.Item(8,7)).Value = $arr
$row = @("wp.pl", "o2.pl")
$arr = New-Object 'object[,]' 1,2
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(9,2), $newws.Cells.Item(9,3)).Value = $arr
$row = @("web.de", "gmx.de", "gmx.ch", "gmx.net", "gmx.com", "gmx.at", "gmx.fr", "mail.com", "1and1.com", "1und1.de")
$arr = New-Object 'object[,]' 1,10
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(10,2), $newws.Cells.Item(10,11)).Value = $arr
$row = @("bigpond.com", "bigpond.net.au", "bigpond.com.au", "telstra.com", "bigpond.net")
$arr = New-Object 'object[,]' 1,5
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(11,2), $newws.Cells.Item(11,6)).Value = $arr
$row = @("softbank.ne.jp", "c.vodafone.ne.jp", "d.vodafone.ne.jp", "h.vodafone.ne.jp", "k.vodafone.ne.jp", "n.vodafone.ne.jp", "q.vodafone.ne.jp", "r.vodafone.ne.jp", "s.vodafone.ne.jp", "t.vodafone.ne.jp", "jp-c.ne.jp", "jp-d.ne.jp", "jp-h.ne.jp", "jp-k.ne.jp", "jp-n.ne.jp", "jp-q.ne.jp", "jp-r.ne.jp", "jp-s.ne.jp", "jp-t.ne.jp")
$arr = New-Object 'object[,]' 1,19
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(13,2), $newws.Cells.Item(13,20)).Value = $arr
$row = @("au.com", "ezweb.ne.jp", "uqmobile.jp")
$arr = New-Object 'object[,]' 1,3
for ($j = 0; $j -lt $row.Count; $j++) { $arr[0,$j] = $row[$j] }
$newws.Range($newws.Cells.Item(14,2), $newws.Cells.Item(14,4)).Value = $arr

# --- Column A width on the new sheet (matches source: autofit-ish width) ---
$newws.Columns.Item(1).ColumnWidth = 13.66

# --- Selections: set each sheet's own cursor position, then land back on "Warmup Plan" ---
$newws.Activate()
[void]$newws.Range("E8").Select()

$ws2.Activate()
[void]$ws2.Range("C6").Select()

$ws1.Activate()
